# Cucumber Class 4 commit 1
# Updates the HRMS testdata workbook:
#  - appends an extra "8" to the three generated usernames
#  - widens the Photograph/Username columns (D, E) to fit the new values
#  - moves the saved cell selection to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeData")

# Username column (E) values for rows 2-4 gained a trailing "8"
$ws.Range("E2").Value = "anasule0012345678"
$ws.Range("E3").Value = "blakenailya0012345678"
$ws.Range("E4").Value = "mikeaj0012345678"

# Widen columns D (Photograph) and E (Username) to fit the longer text
# (target character widths are 40.29 and 27.86; values below are the inputs
# that land the engine's internal pixel-quantized column width closest to
# those targets)
$ws.Columns.Item(4).ColumnWidth = 39.5
$ws.Columns.Item(5).ColumnWidth = 27

# Remember the last selected cell as G11
$ws.Range("G11").Select()

$wb.Save()
